$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# Delete the "Strain" row (row 3) from sheet1 ("openbis-metadata"),
# shifting the remaining rows (Timepoint Type ... Scale) up by one.
$ws1.Rows.Item(3).Delete()

# sheet2 ("openbis-data") A1 changes from "Abs" to "Strain"
$ws2.Range("A1").Value = "Strain"

# Update the recorded cell selections to match the new state.
# Select sheet2's cell first, then sheet1's, so that sheet1 ends up
# as the active sheet (matching the original tabSelected state).
$ws2.Range("A2").Select()
$ws1.Range("A11").Select()
